# Update the LAARI artifact-import template: new column set (A:K instead of
# A:L), new header labels, new sample row, wider columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1) -------------------------------------------------
$headers = @(
    "nome_artefato",
    "codigo_artefato",
    "data_descoberta",
    "tipo",
    "local_origem",
    "localizacao_arqueologica",
    "profundidade",
    "nivel_estratigrafico",
    "coordenadas",
    "estado_conservacao",
    "observacoes"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- New sample row (row 2) --------------------------------------------------
$ws.Cells.Item(2, 1).Value = "Vaso Cerâmico Tupi"
$ws.Cells.Item(2, 2).Value = ""

# "2024-01-15" looks like a date to the smart-parsing Value setter, which
# would silently turn it into a date serial instead of the literal text the
# template expects. Force text mode, write it, then drop back to the plain
# border+wrap formatting shared by the rest of the data row so the cell
# lands on the same cell style as its neighbours instead of a one-off
# text-format style.
$dateCell = $ws.Cells.Item(2, 3)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2024-01-15"
$dateCell.ClearFormats()
$dateCell.Borders.LineStyle = 1
$dateCell.Borders.Weight = 2
$dateCell.WrapText = $true

$ws.Cells.Item(2, 4).Value = "vaso cerâmico"
$ws.Cells.Item(2, 5).Value = "Sítio Arqueológico Exemplo, SP"
$ws.Cells.Item(2, 6).Value = "Setor A, Quadra 5"
$ws.Cells.Item(2, 7).Value = "1.20m"
$ws.Cells.Item(2, 8).Value = "Nível 3"
$ws.Cells.Item(2, 9).Value = "-23.5505, -46.6333"
$ws.Cells.Item(2, 10).Value = "bom"
$ws.Cells.Item(2, 11).Value = "Fragmento de borda com decoração pintada"

# --- Widen the remaining columns (A:K) to raw width 20 ----------------------
# COM's ColumnWidth is expressed in characters of the Normal font and differs
# from the raw OOXML <col width> by the standard ~0.833 padding constant, so
# back that out to land exactly on a stored width of 20.
$ws.Range("A1:K1").EntireColumn.ColumnWidth = 19.166666666666668

# --- Drop the now-unused last column (L) ------------------------------------
$ws.Columns("L").Delete() | Out-Null
